# Atualização automática de preços de eletricidade
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 45905

$ws.Range("B2").Value = 89.25
$ws.Range("C2").Value = 79.02
$ws.Range("D2").Value = 65
$ws.Range("E2").Value = 60
$ws.Range("F2").Value = 53.31
$ws.Range("G2").Value = 60
$ws.Range("H2").Value = 76.79000000000001
$ws.Range("I2").Value = 89.98999999999999
$ws.Range("J2").Value = 88.22
$ws.Range("K2").Value = 79.09999999999999
$ws.Range("L2").Value = 35
$ws.Range("M2").Value = 10
$ws.Range("N2").Value = 3.52
$ws.Range("O2").Value = 3.52
$ws.Range("P2").Value = 3
$ws.Range("Q2").Value = 4.31
$ws.Range("R2").Value = 4.31
$ws.Range("S2").Value = 4.31
$ws.Range("T2").Value = 45.23
$ws.Range("U2").Value = 89.25
$ws.Range("V2").Value = 103.97
$ws.Range("W2").Value = 107.1
$ws.Range("X2").Value = 95
$ws.Range("Y2").Value = 88.38
$ws.Range("Z2").Value = 55.73

$ws.Range("AB2").Value = 98.61
$ws.Range("AD2").Value = 105.54
$ws.Range("AF2").Value = 91.69

$ws.Range("AG2").Value = "4h-18h"
